$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string "TEST" by writing it once so it lands in sharedStrings.xml.
# (written into D2 below, overwriting the copied value)

# 1. Duplicate column C (values + formatting) into column D.
$ws.Range("C1:C24").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C1:C24").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

# 2. Column D differs from the copied column C in two cells.
$ws.Range("D2").Value2 = "TEST"
$ws.Range("D8").Value2 = 3

# 3. Match column D's width/bestFit/style to column C.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# 4. Update the active selection.
$ws.Range("G11").Select()
